# Auto-generated edit script: update Leve profit-calc columns (H-N)
# per scheduled market-data refresh, across ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2834
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 2751
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 2751
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -3403
$ws.Range("H61").Value = 5250
$ws.Range("I61").Value = 9000
$ws.Range("K61").Value = 27000
$ws.Range("M61").Value = -26828
$ws.Range("H97").Value = 4999.6
$ws.Range("I97").Value = 4999
$ws.Range("J97").Value = 4999.75
$ws.Range("K97").Value = 14997
$ws.Range("L97").Value = 14999.25
$ws.Range("M97").Value = -14501
$ws.Range("N97").Value = -15991.25
$ws.Range("H101").Value = 2814.6
$ws.Range("I101").Value = 2796
$ws.Range("K101").Value = 8388
$ws.Range("M101").Value = -6766
$ws.Range("H113").Value = 2238.3845
$ws.Range("I113").Value = 2312.5
$ws.Range("J113").Value = 2205.4443
$ws.Range("K113").Value = 2312.5
$ws.Range("L113").Value = 2205.4443
$ws.Range("M113").Value = 941.5
$ws.Range("N113").Value = -8713.444299999999
$ws.Range("H129").Value = 2223.5
$ws.Range("I129").Value = 2131.3333
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 6393.999899999999
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = -1393.999899999999
$ws.Range("N129").Value = -17500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3486.818
$ws.Range("I45").Value = 1512.3334
$ws.Range("J45").Value = 4227.25
$ws.Range("K45").Value = 1512.3334
$ws.Range("L45").Value = 4227.25
$ws.Range("M45").Value = -1135.3334
$ws.Range("N45").Value = -4981.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2444.3076
$ws.Range("I58").Value = 1893.3334
$ws.Range("K58").Value = 1893.3334
$ws.Range("M58").Value = -1690.3334
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H105").Value = 1486.3334
$ws.Range("I105").Value = 1076.4
$ws.Range("K105").Value = 1076.4
$ws.Range("M105").Value = 670.5999999999999
$ws.Range("H107").Value = 333.23077
$ws.Range("I107").Value = 259.22223
$ws.Range("J107").Value = 499.75
$ws.Range("K107").Value = 259.22223
$ws.Range("L107").Value = 499.75
$ws.Range("M107").Value = 1660.77777
$ws.Range("N107").Value = -4339.75
$ws.Range("H136").Value = 2444.3076
$ws.Range("I136").Value = 1893.3334
$ws.Range("K136").Value = 5680.0002
$ws.Range("M136").Value = -3130.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 751.25
$ws.Range("I97").Value = 668.3333
$ws.Range("K97").Value = 2004.9999
$ws.Range("M97").Value = -1508.9999
$ws.Range("H117").Value = 653.8
$ws.Range("J117").Value = 603
$ws.Range("L117").Value = 1809
$ws.Range("N117").Value = -8693
$ws.Range("H121").Value = 1523.75
$ws.Range("I121").Value = 477.6
$ws.Range("J121").Value = 2271
$ws.Range("K121").Value = 1432.8
$ws.Range("L121").Value = 6813
$ws.Range("M121").Value = -122.8000000000002
$ws.Range("N121").Value = -9433
$ws.Range("H122").Value = 1094.1538
$ws.Range("J122").Value = 1294.6
$ws.Range("L122").Value = 11651.4
$ws.Range("N122").Value = -16551.4
$ws.Range("H125").Value = 17499.834
$ws.Range("I125").Value = 14999.75
$ws.Range("J125").Value = 22500
$ws.Range("K125").Value = 44999.25
$ws.Range("L125").Value = 67500
$ws.Range("M125").Value = -40079.25
$ws.Range("N125").Value = -77340
$ws.Range("H129").Value = 2395.9285
$ws.Range("I129").Value = 1258
$ws.Range("J129").Value = 3249.375
$ws.Range("K129").Value = 3774
$ws.Range("L129").Value = 9748.125
$ws.Range("M129").Value = 1226
$ws.Range("N129").Value = -19748.125
$ws.Range("H131").Value = 372.83334
$ws.Range("I131").Value = 261.75
$ws.Range("K131").Value = 785.25
$ws.Range("M131").Value = 4254.75
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 333.33334
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 1000.00002
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = 4069.99998
$ws.Range("N134").Value = -19140
$ws.Range("H140").Value = 1220.9
$ws.Range("I140").Value = 1023.3333
$ws.Range("K140").Value = 3069.9999
$ws.Range("M140").Value = 2110.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2533.3333
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122
$ws.Range("H126").Value = 3828.5
$ws.Range("I126").Value = 2057.5
$ws.Range("J126").Value = 5599.5
$ws.Range("K126").Value = 6172.5
$ws.Range("L126").Value = 16798.5
$ws.Range("M126").Value = -3702.5
$ws.Range("N126").Value = -21738.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1998.8334
$ws.Range("I7").Value = 1598.8
$ws.Range("J7").Value = 3999
$ws.Range("K7").Value = 1598.8
$ws.Range("L7").Value = 3999
$ws.Range("M7").Value = -1486.8
$ws.Range("N7").Value = -4223
$ws.Range("H43").Value = 29999
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 29999
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 29999
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -30385
$ws.Range("H126").Value = 1998.8334
$ws.Range("I126").Value = 1598.8
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 4796.4
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -2326.4
$ws.Range("N126").Value = -16937
$ws.Range("H136").Value = 3360.8333
$ws.Range("I136").Value = 2833
$ws.Range("K136").Value = 8499
$ws.Range("M136").Value = -5949

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33845.23
$ws.Range("I54").Value = 20000
$ws.Range("J54").Value = 34999
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 34999
$ws.Range("M54").Value = -19480
$ws.Range("N54").Value = -36039
$ws.Range("H103").Value = 41666.668
$ws.Range("J103").Value = 41666.668
$ws.Range("L103").Value = 41666.668
$ws.Range("N103").Value = -44010.668
$ws.Range("H107").Value = 373.2857
$ws.Range("I107").Value = 398.33334
$ws.Range("J107").Value = 354.5
$ws.Range("K107").Value = 1195.00002
$ws.Range("L107").Value = 1063.5
$ws.Range("M107").Value = 724.9999800000001
$ws.Range("N107").Value = -4903.5
